$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (ECs -> FAPs) ---
$ws.Range("G2").Value = 2.119667333333334
$ws.Range("H2").Value = 6.359002
$ws.Range("I2").Value = 0.6371329247828699
$ws.Range("J2").Value = 0.6371329247828699
$ws.Range("M2").Value = 0.06617233333333333
$ws.Range("O2").Value = 0.9596962108540322
$ws.Range("P2").Value = 0.9596962108540322
$ws.Range("Q2").Value = 0.1402633333371111
$ws.Range("R2").Value = 1.262370000034
$ws.Range("S2").Value = 0.6114540537244674
$ws.Range("T2").Value = 0.6114540537244674

# --- Update existing row 3: becomes ECs -> MuSCs ---
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 2.119667333333334
$ws.Range("H3").Value = 6.359002
$ws.Range("I3").Value = 0.6371329247828699
$ws.Range("J3").Value = 0.6371329247828699
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002779
$ws.Range("N3").Value = 0.008337000000000001
$ws.Range("O3").Value = 0.04030378914596769
$ws.Range("P3").Value = 0.04030378914596769
$ws.Range("Q3").Value = 0.005890555519333335
$ws.Range("R3").Value = 0.05301499967400001
$ws.Range("S3").Value = 0.02567887105840248
$ws.Range("T3").Value = 0.02567887105840248

# --- Update existing row 4: becomes FAPs -> FAPs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("G4").Value = 0.8424356666666667
$ws.Range("H4").Value = 2.527307
$ws.Range("I4").Value = 0.2532206312773955
$ws.Range("J4").Value = 0.2532206312773955
$ws.Range("M4").Value = 0.06617233333333333
$ws.Range("O4").Value = 0.9596962108540322
$ws.Range("P4").Value = 0.9596962108540322
$ws.Range("Q4").Value = 0.05574593374655556
$ws.Range("R4").Value = 0.501713403719
$ws.Range("S4").Value = 0.2430148803469825
$ws.Range("T4").Value = 0.2430148803469825

# --- Update existing row 5: becomes FAPs -> MuSCs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.8424356666666667
$ws.Range("H5").Value = 2.527307
$ws.Range("I5").Value = 0.2532206312773955
$ws.Range("J5").Value = 0.2532206312773955
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.002779
$ws.Range("N5").Value = 0.008337000000000001
$ws.Range("O5").Value = 0.04030378914596769
$ws.Range("P5").Value = 0.04030378914596769
$ws.Range("Q5").Value = 0.002341128717666667
$ws.Range("R5").Value = 0.021070158459
$ws.Range("S5").Value = 0.01020575093041298
$ws.Range("T5").Value = 0.01020575093041298

# --- New row 6: MuSCs -> FAPs ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Rln1"
$ws.Range("C6").Value = "Rxfp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2981566666666667
$ws.Range("H6").Value = 0.89447
$ws.Range("I6").Value = 0.08962039754517039
$ws.Range("J6").Value = 0.08962039754517039
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06617233333333333
$ws.Range("N6").Value = 0.198517
$ws.Range("O6").Value = 0.9596962108540322
$ws.Range("P6").Value = 0.9596962108540322
$ws.Range("Q6").Value = 0.01972972233222222
$ws.Range("R6").Value = 0.17756750099
$ws.Range("S6").Value = 0.08600835593933204
$ws.Range("T6").Value = 0.08600835593933204

# --- New row 7: MuSCs -> MuSCs ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Rln1"
$ws.Range("C7").Value = "Rxfp2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2981566666666667
$ws.Range("H7").Value = 0.89447
$ws.Range("I7").Value = 0.08962039754517039
$ws.Range("J7").Value = 0.08962039754517039
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.002779
$ws.Range("N7").Value = 0.008337000000000001
$ws.Range("O7").Value = 0.04030378914596769
$ws.Range("P7").Value = 0.04030378914596769
$ws.Range("Q7").Value = 0.0008285773766666668
$ws.Range("R7").Value = 0.00745719639
$ws.Range("S7").Value = 0.003612041605838348
$ws.Range("T7").Value = 0.003612041605838348

# --- New row 8: Resolving-Mac -> FAPs ---
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Rln1"
$ws.Range("C8").Value = "Rxfp2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.06662433333333333
$ws.Range("H8").Value = 0.199873
$ws.Range("I8").Value = 0.0200260463945642
$ws.Range("J8").Value = 0.0200260463945642
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.06617233333333333
$ws.Range("N8").Value = 0.198517
$ws.Range("O8").Value = 0.9596962108540322
$ws.Range("P8").Value = 0.9596962108540322
$ws.Range("Q8").Value = 0.004408687593444444
$ws.Range("R8").Value = 0.039678188341
$ws.Range("S8").Value = 0.01921892084325031
$ws.Range("T8").Value = 0.01921892084325032

# --- New row 9: Resolving-Mac -> MuSCs ---
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Rln1"
$ws.Range("C9").Value = "Rxfp2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.06662433333333333
$ws.Range("H9").Value = 0.199873
$ws.Range("I9").Value = 0.0200260463945642
$ws.Range("J9").Value = 0.0200260463945642
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.002779
$ws.Range("N9").Value = 0.008337000000000001
$ws.Range("O9").Value = 0.04030378914596769
$ws.Range("P9").Value = 0.04030378914596769
$ws.Range("Q9").Value = 0.0001851490223333333
$ws.Range("R9").Value = 0.001666341201
$ws.Range("S9").Value = 0.0008071255513138818
$ws.Range("T9").Value = 0.0008071255513138819
